$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "unit"
$ws.Range("B1").Value = "question"
$ws.Range("C1").Value = "marks"

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3.1
$ws.Range("A5").Value = 3.2
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

$ws.Columns.Item(1).ColumnWidth = 8.5
$ws.Columns.Item(2).ColumnWidth = 21.67

$ws.Range("B2").Select() | Out-Null
